$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51 (shifts existing rows 51..185 down to 52..186)
$ws.Rows(51).Insert()

# Populate the newly inserted row with the new entry
$ws.Range("A51").Value = "GEP48591"
$ws.Range("B51").Value = "ipb2-EtbP"
